$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price report was inserted. Insert a row at 97 (shifting the
# existing rows 97-124 down to 98-125, which also keeps row 125's data
# identical to the former row 124) and populate the new row with this
# week's "Camote" entry.
$ws.Rows(97).Insert()

$ws.Cells.Item(97, 1).Value = 7
$ws.Cells.Item(97, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(97, 3).Value = "Ñuble"
$ws.Cells.Item(97, 4).Value = 44588
$ws.Cells.Item(97, 5).Value = 16
$ws.Cells.Item(97, 6).Value = 100112045
$ws.Cells.Item(97, 7).Value = "Zapallo"
$ws.Cells.Item(97, 8).Value = "Camote"
$ws.Cells.Item(97, 9).Value = "1a nueva(o)"
$ws.Cells.Item(97, 10).Value = 400
$ws.Cells.Item(97, 11).Value = 300
$ws.Cells.Item(97, 12).Value = 350
$ws.Cells.Item(97, 13).Value = 325
$ws.Cells.Item(97, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(97, 15).Value = "Región del Maule"
$ws.Cells.Item(97, 16).Value = 325
$ws.Cells.Item(97, 17).Value = 1
$ws.Cells.Item(97, 18).Value = "Hortaliza"
